$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.581.76"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.043.55"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'244.99"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'55.74"
$ws.Range("E8").Value = "  -3.42%  "
$ws.Range("D9").Value = "'63.49"
$ws.Range("E9").Value = "  +7.36%  "
$ws.Range("D10").Value = "'0.367"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "'0.0747"
$ws.Range("E11").Value = "  -3.89%  "
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").Value = "'0.912"
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("D14").Value = "'14.51"
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("D15").Value = "2.342.76"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "'5.40"
$ws.Range("E16").Value = "  -2.96%  "
$ws.Range("D17").Value = "2.072.16"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "36.448.11"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "'17.26"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "'72.01"
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("D21").Value = "0.0₃0855"
$ws.Range("E21").Value = "  -3.82%  "
$ws.Range("D22").Value = "'237.83"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").Value = "'5.17"
$ws.Range("E23").Value = "  -4.63%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("D27").Value = "'9.19"
$ws.Range("E27").Value = "  -8.44%  "
$ws.Range("D28").Value = "'164.28"
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("D29").Value = "'19.98"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").Value = "'0.121"
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("E32").Value = "  -8.34%  "
$ws.Range("D33").Value = "'0.0599"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("D34").Value = "'4.42"
$ws.Range("E34").Value = "  -7.75%  "
$ws.Range("D35").Value = "'0.0871"
$ws.Range("E35").Value = "  +3.20%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("D38").Value = "'2.19"
$ws.Range("E38").Value = "  -8.69%  "
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").Value = "'4.99"
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.23"
$ws.Range("E40").Value = "  -5.96%  "
$ws.Range("E41").Value = "  -3.40%  "
$ws.Range("E42").Value = "  -3.44%  "
$ws.Range("D43").Value = "'1.10"
$ws.Range("E43").Value = "  -3.78%  "
$ws.Range("D44").Value = "'93.37"
$ws.Range("E44").Value = "  -3.61%  "
$ws.Range("D45").Value = "'0.0903"
$ws.Range("E45").Value = "  -6.23%  "
$ws.Range("D46").Value = "'15.89"
$ws.Range("E46").Value = "  -5.22%  "
$ws.Range("D47").Value = "'7.52"
$ws.Range("E47").Value = "  +10.94%  "
$ws.Range("D48").Value = "1.374.29"
$ws.Range("E48").Value = "  +4.66%  "
$ws.Range("D49").Value = "'2.93"
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("D50").Value = "'2.25"
$ws.Range("E50").Value = "  -5.30%  "
$ws.Range("D51").Value = "'45.70"
$ws.Range("E51").Value = "  +0.90%  "
